$wb = $excel.ActiveWorkbook

# --- Hoja1: update the "Conversión del día" note with the refreshed rates ---
$wsHoja1 = $wb.Worksheets.Item("Hoja1")
$newText = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 7.02 = 28117.12 pesos`n✅ 28117.12 pesos = 7.0 = 974.4 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"
$wsHoja1.Range("A1").Value = $newText

# --- tasas: refresh the rate table cells ---
$wsTasas = $wb.Worksheets.Item("tasas")
$wsTasas.Range("N10").Value = 142.5
$wsTasas.Range("O10").Value = 4006.69
$wsTasas.Range("N12").Value = 4017
$wsTasas.Range("O12").Value = 139.21
